# Applies the "final results" data refresh described in the commit:
#   - CAPEX!W2/X2 bumped to the new placeholder totals
#   - CAPEX!S11/S12/S15 (Duct Cost) re-shuffled between the two cached figures
#   - CAPEX!U3:U15 / V3:V15 (CO cost / RN cost) updated to the re-run figures
#   - OPEX!B2:B15 frozen from formulas into their latest computed values
#   - Selections / active sheet restored to match the saved-state of the workbook

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# CAPEX sheet
# ---------------------------------------------------------------------------
$capex = $wb.Worksheets.Item("CAPEX")

$capex.Range("W2").Value = 100000
$capex.Range("X2").Value = 100000

$capex.Range("U3").Value = 3419.7777777777778
$capex.Range("V3").Value = 119020.4

$capex.Range("U4").Value = 7900
$capex.Range("V4").Value = 116275.2

$capex.Range("U5").Value = 7232.1111111111113
$capex.Range("V5").Value = 98960

$capex.Range("U6").Value = 8289
$capex.Range("V6").Value = 341360

$capex.Range("U7").Value = 30450
$capex.Range("V7").Value = 327450

$capex.Range("U8").Value = 19660.666666666668
$capex.Range("V8").Value = 166336.79999999999

$capex.Range("U9").Value = 17040
$capex.Range("V9").Value = 195632.2

$capex.Range("U10").Value = 15506
$capex.Range("V10").Value = 201120

$capex.Range("S11").Value = 148601.5257512136
$capex.Range("U11").Value = 5360
$capex.Range("V11").Value = 79394

$capex.Range("S12").Value = 162896.27914347179
$capex.Range("U12").Value = 12380
$capex.Range("V12").Value = 89994

$capex.Range("U13").Value = 24400
$capex.Range("V13").Value = 228180

$capex.Range("U14").Value = 38000
$capex.Range("V14").Value = 272276

$capex.Range("S15").Value = 162896.27914347179
$capex.Range("U15").Value = 22000
$capex.Range("V15").Value = 215380

# ---------------------------------------------------------------------------
# OPEX sheet - the AVERAGE()/1.5*x formulas are frozen to their latest values
# ---------------------------------------------------------------------------
$opex = $wb.Worksheets.Item("OPEX")

$opex.Range("B2").Value = 10000
$opex.Range("B3").Value = 42365.673348618162
$opex.Range("B4").Value = 42358.112806224512
$opex.Range("B5").Value = 42902.769000487555
$opex.Range("B6").Value = 35408.524500181637
$opex.Range("B7").Value = 40576.763975644033
$opex.Range("B8").Value = 37982.807078103164
$opex.Range("B9").Value = 43544.075083463555
$opex.Range("B10").Value = 45178.949649365633
$opex.Range("B11").Value = 22360.217706464937
$opex.Range("B12").Value = 24572.597341958404
$opex.Range("B13").Value = 43083.217137892461
$opex.Range("B14").Value = 52368.522953392559
$opex.Range("B15").Value = 49441.019044068453

# ---------------------------------------------------------------------------
# Restore cursor/selection state seen in the saved workbook
# ---------------------------------------------------------------------------
$capex.Activate()
$capex.Range("N29").Select()

$opex.Activate()
$opex.Range("B26").Select()
